$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 253, shifting existing rows 253:301 down to 254:302.
$ws.Rows("253:253").Insert()

# Populate the newly inserted row 253 with the new record.
$ws.Range("A253").Value = 11
$ws.Range("B253").Value = "Vega Monumental Concepción"
$ws.Range("C253").Value = "Bíobío"
$ws.Range("D253").Value = 44785
$ws.Range("E253").Value = 8
$ws.Range("F253").Value = 100114001
$ws.Range("G253").Value = "Papa"
$ws.Range("H253").Value = "Asterix"
$ws.Range("I253").Value = "1a (guarda)"
$ws.Range("J253").Value = 4000
$ws.Range("K253").Value = 8000
$ws.Range("L253").Value = 8500
$ws.Range("M253").Value = 8250
$ws.Range("N253").Value = "$/saco 25 kilos"
$ws.Range("O253").Value = "Provincia de Arauco"
$ws.Range("P253").Value = 330
$ws.Range("Q253").Value = 25
$ws.Range("R253").Value = "Hortaliza"
